$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestCases")
$ws2 = $wb.Worksheets.Item("TestData")

# --- Runmode value flips (common utility for runmodes) ---
# TestCases: OpenAccountTest Runmode N -> Y
$ws1.Range("B3").Value = "Y"

# TestData: AddCustomerTest second data row Runmode N -> Y
$ws2.Range("A4").Value = "Y"
# TestData: OpenAccountTest data rows Runmode Y -> N
$ws2.Range("A8").Value = "N"
$ws2.Range("A9").Value = "N"

# --- Selection / active sheet state ---
# TestData sheet: selection moves to A4, sheet no longer the active tab
$ws2.Activate()
$ws2.Range("A4").Select()

# TestCases sheet becomes the active tab, selection at B4
$ws1.Activate()
$ws1.Range("B4").Select()
